# Update the Liph-Lpar1 LR-pair sheet with newly recomputed TPM-based values.
# The underlying recompute (a) dropped the "MuSCs" sending-cluster block
# (rows that used to be 10-13) and (b) refreshed every specificity /
# weight metric for the remaining ECs- and FAPs-sending rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows for Sending cluster = MuSCs (old rows 10-13).
$ws.Rows("10:13").Delete()

$data = @(
    @("ECs", "Liph", "Lpar1", "ECs", 1, 0.3333333333333333, 0.04577266666666666, 0.137318, 0.6974138631561839, 0.697413863156184, 2, 0.6666666666666666, 0.3268783333333333, 0.9806349999999999, 0.003643071043709958, 0.003643071043709957, 0.01496209299222222, 0.13465883693, 0.002540728250346192, 0.002540728250346193),
    @("ECs", "Liph", "Lpar1", "FAPs", 1, 0.3333333333333333, 0.04577266666666666, 0.137318, 0.6974138631561839, 0.697413863156184, 3, 1, 73.66398133333333, 220.991944, 0.820987780447947, 0.820987780447947, 3.371796862910222, 30.346171766192, 0.5725682595662237, 0.5725682595662237),
    @("ECs", "Liph", "Lpar1", "MuSCs", 1, 0.3333333333333333, 0.04577266666666666, 0.137318, 0.6974138631561839, 0.697413863156184, 3, 1, 15.26610033333333, 45.798301, 0.1701412495211907, 0.1701412495211906, 0.6987701218575555, 6.288931096717999, 0.1186588661107938, 0.1186588661107938),
    @("ECs", "Liph", "Lpar1", "Resolving-Mac", 1, 0.3333333333333333, 0.04577266666666666, 0.137318, 0.6974138631561839, 0.697413863156184, 2, 0.6666666666666666, 0.4690786666666667, 1.407236, 0.005227898987152434, 0.005227898987152434, 0.02147098144977778, 0.193238833048, 0.00364600922882028, 0.00364600922882028),
    @("FAPs", "Liph", "Lpar1", "ECs", 1, 0.3333333333333333, 0.01985933333333333, 0.059578, 0.302586136843816, 0.302586136843816, 2, 0.6666666666666666, 0.3268783333333333, 0.9806349999999999, 0.003643071043709958, 0.003643071043709957, 0.006491585781111111, 0.05842427202999999, 0.001102342793363765, 0.001102342793363765),
    @("FAPs", "Liph", "Lpar1", "FAPs", 1, 0.3333333333333333, 0.01985933333333333, 0.059578, 0.302586136843816, 0.302586136843816, 3, 1, 73.66398133333333, 220.991944, 0.820987780447947, 0.820987780447947, 1.462917559959111, 13.166258039632, 0.2484195208817233, 0.2484195208817233),
    @("FAPs", "Liph", "Lpar1", "MuSCs", 1, 0.3333333333333333, 0.01985933333333333, 0.059578, 0.302586136843816, 0.302586136843816, 3, 1, 15.26610033333333, 45.798301, 0.1701412495211907, 0.1701412495211906, 0.3031745752197778, 2.728571176977999, 0.05148238341039685, 0.05148238341039684),
    @("FAPs", "Liph", "Lpar1", "Resolving-Mac", 1, 0.3333333333333333, 0.01985933333333333, 0.059578, 0.302586136843816, 0.302586136843816, 2, 0.6666666666666666, 0.4690786666666667, 1.407236, 0.005227898987152434, 0.005227898987152434, 0.00931558960088889, 0.08384030640800001, 0.001581889758332154, 0.001581889758332154)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
